# newAssignment return from server
# Update shift start/end times (columns C/D) and duration-derived value (column E)
# to reflect the newly returned assignment schedule from the server.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "08:00"
$ws.Range("D1").Value = "13:00"
$ws.Range("C2").Value = "12:00"
$ws.Range("D2").Value = "17:00"
$ws.Range("D3").Value = "14:00"
$ws.Range("E3").Value = 60
$ws.Range("C4").Value = "09:00"
$ws.Range("D4").Value = "14:00"
$ws.Range("C5").Value = "08:00"
$ws.Range("D5").Value = "11:00"
$ws.Range("C6").Value = "14:00"
$ws.Range("E6").Value = 60
$ws.Range("D7").Value = "13:00"
$ws.Range("E7").Value = 60
$ws.Range("C8").Value = "07:00"
$ws.Range("D8").Value = "12:00"
$ws.Range("E8").Value = 80
$ws.Range("C9").Value = "11:00"
$ws.Range("D9").Value = "14:00"
$ws.Range("C10").Value = "08:00"
$ws.Range("D10").Value = "11:00"
$ws.Range("C11").Value = "09:00"
$ws.Range("D11").Value = "12:00"
$ws.Range("E11").Value = 60
$ws.Range("D12").Value = "19:00"
$ws.Range("E12").Value = 80
$ws.Range("C13").Value = "08:00"
$ws.Range("D13").Value = "11:00"
$ws.Range("C14").Value = "15:00"
$ws.Range("D14").Value = "18:00"
$ws.Range("C15").Value = "10:00"
$ws.Range("D15").Value = "13:00"
$ws.Range("D18").Value = "19:00"
$ws.Range("E18").Value = 80
$ws.Range("C19").Value = "07:00"
$ws.Range("D19").Value = "12:00"
$ws.Range("C20").Value = "11:00"
$ws.Range("D20").Value = "14:00"
$ws.Range("E20").Value = 60
$ws.Range("D21").Value = "14:00"
$ws.Range("E21").Value = 60
$ws.Range("D22").Value = "15:00"
$ws.Range("E22").Value = 80
$ws.Range("D23").Value = "12:00"
$ws.Range("E23").Value = 60
$ws.Range("C24").Value = "07:00"
$ws.Range("D24").Value = "12:00"
$ws.Range("E24").Value = 80
$ws.Range("C26").Value = "08:00"
$ws.Range("D26").Value = "11:00"
$ws.Range("C27").Value = "09:00"
$ws.Range("D27").Value = "14:00"
$ws.Range("C28").Value = "08:00"
$ws.Range("D28").Value = "11:00"
$ws.Range("E28").Value = 60
$ws.Range("C29").Value = "07:00"
$ws.Range("D29").Value = "10:00"
$ws.Range("C30").Value = "14:00"
$ws.Range("D30").Value = "19:00"
$ws.Range("E30").Value = 80
$ws.Range("C31").Value = "10:00"
$ws.Range("D31").Value = "13:00"
$ws.Range("C32").Value = "13:00"
$ws.Range("D32").Value = "18:00"
$ws.Range("C33").Value = "07:00"
$ws.Range("D33").Value = "10:00"
$ws.Range("E33").Value = 60
$ws.Range("D34").Value = "13:00"
$ws.Range("E34").Value = 80
$ws.Range("C35").Value = "07:00"
$ws.Range("D35").Value = "12:00"
$ws.Range("E35").Value = 80
$ws.Range("C36").Value = "15:00"
$ws.Range("D36").Value = "18:00"
$ws.Range("C37").Value = "09:00"
$ws.Range("D37").Value = "14:00"
$ws.Range("C38").Value = "10:00"
$ws.Range("D38").Value = "13:00"
$ws.Range("E38").Value = 60
$ws.Range("C40").Value = "10:00"
$ws.Range("D40").Value = "15:00"
$ws.Range("D41").Value = "16:00"
$ws.Range("E41").Value = 60
$ws.Range("C42").Value = "10:00"
$ws.Range("E42").Value = 60
$ws.Range("C43").Value = "12:00"
$ws.Range("D43").Value = "17:00"
$ws.Range("C44").Value = "12:00"
$ws.Range("D44").Value = "15:00"
$ws.Range("D45").Value = "15:00"
$ws.Range("E45").Value = 60
$ws.Range("C46").Value = "10:00"
$ws.Range("D46").Value = "15:00"
$ws.Range("C47").Value = "10:00"
$ws.Range("D47").Value = "13:00"
$ws.Range("E47").Value = 60
$ws.Range("C48").Value = "12:00"
$ws.Range("D48").Value = "17:00"
$ws.Range("E48").Value = 80
$ws.Range("C49").Value = "08:00"
$ws.Range("D49").Value = "11:00"
$ws.Range("C50").Value = "08:00"
$ws.Range("D50").Value = "11:00"
$ws.Range("C51").Value = "07:00"
$ws.Range("D51").Value = "12:00"
$ws.Range("C52").Value = "10:00"
$ws.Range("D52").Value = "13:00"
$ws.Range("C53").Value = "14:00"
$ws.Range("D53").Value = "19:00"
$ws.Range("D54").Value = "10:00"
$ws.Range("E54").Value = 60
$ws.Range("C55").Value = "10:00"
$ws.Range("D55").Value = "15:00"
$ws.Range("C56").Value = "14:00"
$ws.Range("D56").Value = "19:00"
$ws.Range("E56").Value = 80
$ws.Range("C57").Value = "10:00"
$ws.Range("D57").Value = "15:00"
$ws.Range("E57").Value = 80
$ws.Range("C58").Value = "11:00"
$ws.Range("D58").Value = "14:00"
$ws.Range("E58").Value = 60
$ws.Range("C59").Value = "07:00"
$ws.Range("D59").Value = "12:00"
$ws.Range("C60").Value = "07:00"
$ws.Range("D60").Value = "10:00"
$ws.Range("C61").Value = "07:00"
$ws.Range("D61").Value = "10:00"
$ws.Range("C62").Value = "08:00"
$ws.Range("D62").Value = "13:00"
$ws.Range("E62").Value = 80
$ws.Range("C63").Value = "13:00"
$ws.Range("D63").Value = "18:00"
$ws.Range("C64").Value = "09:00"
$ws.Range("D64").Value = "14:00"
$ws.Range("E64").Value = 80
$ws.Range("C65").Value = "07:00"
$ws.Range("D65").Value = "12:00"
$ws.Range("E65").Value = 80
$ws.Range("C67").Value = "09:00"
$ws.Range("D67").Value = "14:00"
$ws.Range("E67").Value = 80
$ws.Range("C68").Value = "13:00"
$ws.Range("D68").Value = "18:00"
$ws.Range("E68").Value = 80
$ws.Range("C69").Value = "11:00"
$ws.Range("D69").Value = "14:00"
$ws.Range("C70").Value = "09:00"
$ws.Range("D70").Value = "12:00"
$ws.Range("E70").Value = 60
$ws.Range("C71").Value = "13:00"
$ws.Range("D71").Value = "18:00"
$ws.Range("C72").Value = "10:00"
$ws.Range("D72").Value = "13:00"
$ws.Range("C73").Value = "08:00"
$ws.Range("D73").Value = "11:00"
$ws.Range("E73").Value = 60
$ws.Range("C74").Value = "12:00"
$ws.Range("D74").Value = "17:00"
$ws.Range("E74").Value = 80
$ws.Range("D75").Value = "15:00"
$ws.Range("E75").Value = 60
$ws.Range("C76").Value = "10:00"
$ws.Range("D76").Value = "15:00"
$ws.Range("E76").Value = 80
$ws.Range("C77").Value = "07:00"
$ws.Range("D77").Value = "10:00"
$ws.Range("C78").Value = "10:00"
$ws.Range("D78").Value = "15:00"
$ws.Range("C79").Value = "14:00"
$ws.Range("D79").Value = "19:00"
$ws.Range("E79").Value = 80
$ws.Range("C80").Value = "09:00"
$ws.Range("D80").Value = "12:00"
$ws.Range("E80").Value = 60
$ws.Range("C81").Value = "15:00"
$ws.Range("D81").Value = "20:00"
$ws.Range("E81").Value = 80
$ws.Range("C82").Value = "09:00"
$ws.Range("D82").Value = "14:00"
$ws.Range("E82").Value = 80
$ws.Range("C83").Value = "10:00"
$ws.Range("D83").Value = "13:00"
$ws.Range("E83").Value = 60
$ws.Range("C84").Value = "11:00"
$ws.Range("D84").Value = "16:00"
$ws.Range("E84").Value = 80
$ws.Range("C85").Value = "13:00"
$ws.Range("D85").Value = "16:00"
$ws.Range("E85").Value = 60
$ws.Range("C86").Value = "11:00"
$ws.Range("D86").Value = "16:00"
$ws.Range("C87").Value = "14:00"
$ws.Range("D87").Value = "17:00"
$ws.Range("C88").Value = "10:00"
$ws.Range("D88").Value = "13:00"
$ws.Range("C89").Value = "08:00"
$ws.Range("D89").Value = "13:00"
$ws.Range("E89").Value = 80
$ws.Range("C90").Value = "10:00"
$ws.Range("D90").Value = "15:00"
$ws.Range("E90").Value = 80
$ws.Range("C91").Value = "11:00"
$ws.Range("D91").Value = "14:00"
$ws.Range("C92").Value = "09:00"
$ws.Range("D92").Value = "14:00"
$ws.Range("C93").Value = "08:00"
$ws.Range("D93").Value = "11:00"
$ws.Range("E93").Value = 60
$ws.Range("C94").Value = "08:00"
$ws.Range("D94").Value = "13:00"
$ws.Range("E94").Value = 80
$ws.Range("C95").Value = "12:00"
$ws.Range("D95").Value = "17:00"
